$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = "D6"
$ws.Range("C3").Value = "A1"
$ws.Range("D3").Value = "B2"
$ws.Range("E3").Value = "B3"
$ws.Range("F3").Value = "C4"
$ws.Range("G3").Value = "D5"
Write-Host $ws.Range("B3").Value2
Write-Host $ws.Range("C3").Value2
Write-Host $ws.Range("D3").Value2
Write-Host $ws.Range("E3").Value2
Write-Host $ws.Range("F3").Value2
Write-Host $ws.Range("G3").Value2
